$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns I (I0) and J (IF)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style (bold, centered, bordered) from H1 onto I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for I2:J83
$data = @(
    @(8, 8),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(7, 7),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(8, 9),
    @(7, 7),
    @(9, 9),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(7, 8),
    @(9, 9),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(8, 9),
    @(9, 9),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(7, 8),
    @(9, 9),
    @(7, 7),
    @(8, 8),
    @(9, 9),
    @(7, 8),
    @(8, 8),
    @(8, 8),
    @(7, 7),
    @(7, 8),
    @(7, 8),
    @(7, 7),
    @(8, 8),
    @(7, 8),
    @(9, 9),
    @(8, 8),
    @(7, 8),
    @(6, 7),
    @(7, 7),
    @(8, 9),
    @(7, 8),
    @(6, 6),
    @(8, 8),
    @(8, 8),
    @(7, 8),
    @(8, 8),
    @(8, 9),
    @(8, 8),
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(6, 7),
    @(8, 8),
    @(7, 7),
    @(7, 8),
    @(8, 8),
    @(9, 9),
    @(6, 6),
    @(8, 8),
    @(6, 6),
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(9, 9)
)

for ($k = 0; $k -lt $data.Length; $k++) {
    $row = $k + 2
    $ws.Cells.Item($row, 9).Value = $data[$k][0]
    $ws.Cells.Item($row, 10).Value = $data[$k][1]
}

